# Populate the "credits" worksheet with track/artist credit rows.
# (fix: finalização dos ajustes do storage manager)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("credits")

# ---- Header row (row 1), columns B..K ------------------------------------
$headers = @(
    "trackUri",
    "trackTitle",
    "roleTitle",
    "artist_uri",
    "artist_name",
    "artist_image_uri",
    "artist_subroles",
    "artist_weight",
    "artist_external_url",
    "artist_creator_uri"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

# ---- Data rows (rows 2..13) ------------------------------------------------
# Columns: idx, trackUri, trackTitle, roleTitle, artist_uri, artist_name,
#          artist_image_uri, artist_subroles, artist_weight,
#          artist_external_url, artist_creator_uri
$rows = @(
    @(0,  "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Performers", "spotify:artist:6vWDO969PvNqNYHIOW5v0m", "Beyoncé",                     "https://i.scdn.co/image/ab677762000078e645c984e8c82f9ce15ebf1f51", "['main artist']", 0.8999000191688538, "", ""),
    @(1,  "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Writers",    "",                                        'Brittany "@Chi_Coney" Coney', "",                                                                   "[]",               0,                  "https://artists.spotify.com/songwriter/1uYXu9L0vzo1lcJk7SjgCT", "spotify:songwriter:1uYXu9L0vzo1lcJk7SjgCT"),
    @(2,  "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Writers",    "",                                        'Denisia "Blu June" Andrews',  "",                                                                   "[]",               0,                  "https://artists.spotify.com/songwriter/6usZxRd8e7TQwcksnohR1U", "spotify:songwriter:6usZxRd8e7TQwcksnohR1U"),
    @(3,  "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Writers",    "",                                        "A. Graham",                    "",                                                                   "[]",               0,                  "", ""),
    @(20, "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Producers",  "spotify:artist:2pSxft78wxISNLTbjFThTk",  "Cadenza",                       "https://i.scdn.co/image/ab677762000078e6afc079cda32d54850e82c385", "['producer']",     0.548799991607666,  "", ""),
    @(21, "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Producers",  "spotify:artist:6PIl48ctj6HDmqVM2uIikd",  "Kelman Duran",                  "https://i.scdn.co/image/ab677762000078e645c984e8c82f9ce15ebf1f51", "['producer']",     0.5483999848365784, "", ""),
    @(22, "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Producers",  "spotify:artist:7rzPhcWW1k8njDS5r2xbWe",  "Jahaan Sweet",                  "https://i.scdn.co/image/ab677762000078e6afc079cda32d54850e82c385", "['producer']",     0.5480999946594238, "", ""),
    @(23, "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Producers",  "spotify:artist:2TRy6nJhEnznYDL0AtYJAI",  "Stuart White",                  "https://i.scdn.co/image/ab677762000078e6d2b377637d9f6ed34f1652e2", "['producer']",     0.5473999977111816, "", ""),
    @(24, "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Producers",  "spotify:artist:0NPWbgGDaiX1yZaN0nzhqn",  "Neenyo",                        "https://i.scdn.co/image/ab677762000078e6d2b377637d9f6ed34f1652e2", "['producer']",     0.5472000241279602, "", ""),
    @(25, "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Producers",  "spotify:artist:6T9dsSIQ31Sk8x24lwfTG4",  "Sevn Thomas",                   "https://i.scdn.co/image/ab677762000078e6d2b377637d9f6ed34f1652e2", "['producer']",     0.5468999743461609, "", ""),
    @(26, "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Producers",  "spotify:artist:2RQnTQM6OQd8dTUDF7UsT7",  "Boi-1da",                       "https://i.scdn.co/image/ab677762000078e6aae517f57fa4e9833b4bce9c", "['producer']",     0.5468000173568726, "", ""),
    @(27, "spotify:track:1w7cgGZR86yWz1pA2puVJD", "HEATED", "Producers",  "spotify:artist:3FdxqLD6H5pfO0C0kd9YYt",  "Harry Edwards",                 "https://i.scdn.co/image/ab677762000078e645c984e8c82f9ce15ebf1f51", "['producer']",     0.5461000204086304, "", "")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $rowData = $rows[$r]

    # Column A: numeric row index, bold/centered/bordered like the header style
    $aCell = $ws.Cells.Item($rowNum, 1)
    $aCell.Value = $rowData[0]
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108   # xlCenter
    $aCell.VerticalAlignment = -4160     # xlTop
    $aCell.Borders.LineStyle = 1

    # Columns B..K: plain values (strings / numbers), skip blanks (no cell)
    for ($c = 1; $c -lt $rowData.Length; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val -and $val -ne "") {
            $ws.Cells.Item($rowNum, $c + 1).Value = $val
        }
    }
}

Write-Output "credits sheet populated"
